# Handback report generation: add the new file
# "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md" as a new row (row 4) to every
# sheet's table: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$HL_COLOR = 15570276   # BGR(..) value that maps to RGB FF6495ED (workbook's HyperLink font colour)
$DT_FMT = "yyyy-mm-dd HH:mm:ss"

$SRC_REPO_URL = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11772e054ad11acbf9733480d64c935eba11add6/e2e/c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$ZHCN_REPO_URL = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2a3e135cac06dddb0dcfd2af696fa5ae50236bc/e2e/c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$DEDE_REPO_URL = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2b8ed04e9056f1d7ef59ff0c86a7642890bd625b/e2e/c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $HL_COLOR
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)
$ovTable.ListRows.Add() | Out-Null

$ov.Range("A4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$ov.Range("B4").Value = "e2e\c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$ov.Range("C4").Value = ".md"
$ov.Range("E4").Value = "Handed back: in sync with en-US"
$ov.Range("F4").Value = "Handed back: in sync with en-US"
$ov.Range("G4").Value = "2016-10-19 11:48:52"
$ov.Range("G4").NumberFormat = $DT_FMT

$ov.Hyperlinks.Add($ov.Range("B4"), $SRC_REPO_URL, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md") | Out-Null
Style-AsHyperlink $ov.Range("B4")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)
$zhTable.ListRows.Add() | Out-Null

$zh.Range("A4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("D4").Value = "e2e"
$zh.Range("E4").Value = "ht"
$zh.Range("F4").Value = "True"
$zh.Range("G4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.047c4204302ad5a720637b9fff9542ac43425c1d.zh-cn.xlf"
$zh.Range("H4").Value = "2016-10-19 11:48:41"
$zh.Range("H4").NumberFormat = $DT_FMT
$zh.Range("I4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$zh.Range("J4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.047c4204302ad5a720637b9fff9542ac43425c1d.zh-cn.xlf"
$zh.Range("K4").Value = "2016-10-19 11:49:24"
$zh.Range("K4").NumberFormat = $DT_FMT
$zh.Range("M4").Value = "True"
$zh.Range("O4").Value = "False"

$zh.Hyperlinks.Add($zh.Range("A4"), $SRC_REPO_URL, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md") | Out-Null
Style-AsHyperlink $zh.Range("A4")

$zh.Hyperlinks.Add($zh.Range("I4"), $ZHCN_REPO_URL, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md") | Out-Null
Style-AsHyperlink $zh.Range("I4")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)
$deTable.ListRows.Add() | Out-Null

$de.Range("A4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("D4").Value = "e2e"
$de.Range("E4").Value = "ht"
$de.Range("F4").Value = "True"
$de.Range("G4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.047c4204302ad5a720637b9fff9542ac43425c1d.de-de.xlf"
$de.Range("H4").Value = "2016-10-19 11:48:52"
$de.Range("H4").NumberFormat = $DT_FMT
$de.Range("I4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md"
$de.Range("J4").Value = "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.047c4204302ad5a720637b9fff9542ac43425c1d.de-de.xlf"
$de.Range("K4").Value = "2016-10-19 11:49:42"
$de.Range("K4").NumberFormat = $DT_FMT
$de.Range("M4").Value = "True"
$de.Range("O4").Value = "False"

$de.Hyperlinks.Add($de.Range("A4"), $SRC_REPO_URL, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md") | Out-Null
Style-AsHyperlink $de.Range("A4")

$de.Hyperlinks.Add($de.Range("I4"), $DEDE_REPO_URL, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "c19fa120-f7c3-4788-b67d-b4fcf63aa51f.md") | Out-Null
Style-AsHyperlink $de.Range("I4")

Write-Host "Handback row added to Overview, zh-cn, de-de."
